$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '20.373.85'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -7.47%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.444.20'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -6.98%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.007'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.49%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.003'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '277.71'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3735'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.77%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3070'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '40.52'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -8.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.010'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06556'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.008'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.383'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.22'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -7.46%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.139'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.443.97'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001011'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.53%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.05877'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -10.44%  '
$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '76.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -8.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.737'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -6.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.324'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '20.383.50'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '143.50'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.99%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.219'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -7.11%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.609.40'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '109.60'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.9083'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.631'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -25.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.403'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07743'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.319'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -8.36%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.90'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.07%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.004'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.31%  '
$ws.Range('B38').Value = 'WEMIXTOKEN'
$ws.Range('C38').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.406'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -12.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05633'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.94%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.140'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.28%  '
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.737'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.75%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1913'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.63%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.02037'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.587'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5321'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.06'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.72%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5154'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '111.63'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.783'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.056'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.007'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.52%  '
